$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the MODS wrapper element text in C1 and Z1 to replace
# <update type="MODS"> with <datastream type="md_descriptive" operation="update">
$ws.Range("C1").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink">'
$ws.Range("Z1").Value = "</mods:mods></datastream></object>"

# Reflect the updated selection/view state (C1 selected, no frozen/scrolled topLeftCell)
$ws.Range("C1").Select()
